# Scheduled-runner refresh of cached market-board figures
# (currentAveragePrice / NQ / HQ / LevePrice / LeveProfit columns H:N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 827.5833
$ws.Range("I9").Value = 1015.875
$ws.Range("J9").Value = 451
$ws.Range("K9").Value = 1015.875
$ws.Range("L9").Value = 451
$ws.Range("M9").Value = -846.875
$ws.Range("N9").Value = -789

$ws.Range("H99").Value = 3037.1428
$ws.Range("I99").Value = 4980
$ws.Range("J99").Value = 446.66666
$ws.Range("K99").Value = 14940
$ws.Range("L99").Value = 1339.99998
$ws.Range("M99").Value = -13442
$ws.Range("N99").Value = -4335.999980000001

$ws.Range("H111").Value = 29554.5
$ws.Range("I111").Value = 14444.875
$ws.Range("J111").Value = 59773.75
$ws.Range("K111").Value = 43334.625
$ws.Range("L111").Value = 179321.25
$ws.Range("M111").Value = -40267.625
$ws.Range("N111").Value = -185455.25

$ws.Range("H131").Value = 9288.375
$ws.Range("I131").Value = 3457.1
$ws.Range("J131").Value = 13453.571
$ws.Range("K131").Value = 10371.3
$ws.Range("L131").Value = 40360.713
$ws.Range("M131").Value = -5331.299999999999
$ws.Range("N131").Value = -50440.713

$ws.Range("H132").Value = 2107.5657
$ws.Range("I132").Value = 1831.0946
$ws.Range("J132").Value = 2925.92
$ws.Range("K132").Value = 5493.283799999999
$ws.Range("L132").Value = 8777.76
$ws.Range("M132").Value = -2963.283799999999
$ws.Range("N132").Value = -13837.76

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1002717.4
$ws.Range("I61").Value = 834945.7
$ws.Range("K61").Value = 834945.7
$ws.Range("M61").Value = -834733.7

$ws.Range("H132").Value = 466441.12
$ws.Range("I132").Value = 278957.12
$ws.Range("K132").Value = 836871.36
$ws.Range("M132").Value = -834341.36

$ws.Range("H136").Value = 1002717.4
$ws.Range("I136").Value = 834945.7
$ws.Range("K136").Value = 2504837.1
$ws.Range("M136").Value = -2502287.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 11537.111
$ws.Range("J64").Value = 17056.5
$ws.Range("L64").Value = 17056.5
$ws.Range("N64").Value = -17506.5

$ws.Range("H67").Value = 11537.111
$ws.Range("J67").Value = 17056.5
$ws.Range("L67").Value = 17056.5
$ws.Range("N67").Value = -18616.5

$ws.Range("H94").Value = 3472.5789
$ws.Range("I94").Value = 1607.3
$ws.Range("J94").Value = 5545.1113
$ws.Range("K94").Value = 1607.3
$ws.Range("L94").Value = 5545.1113
$ws.Range("M94").Value = -1156.3
$ws.Range("N94").Value = -6447.1113

$ws.Range("H107").Value = 9624727
$ws.Range("I107").Value = 6628.8945
$ws.Range("K107").Value = 6628.8945
$ws.Range("M107").Value = -4708.8945

$ws.Range("H134").Value = 1691.8667
$ws.Range("I134").Value = 1542.28
$ws.Range("K134").Value = 4626.84
$ws.Range("M134").Value = -2091.84

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2703.965
$ws.Range("I31").Value = 1479.8572
$ws.Range("J31").Value = 2941.986
$ws.Range("K31").Value = 1479.8572
$ws.Range("L31").Value = 2941.986
$ws.Range("M31").Value = -1184.8572
$ws.Range("N31").Value = -3531.986

$ws.Range("H34").Value = 2703.965
$ws.Range("I34").Value = 1479.8572
$ws.Range("J34").Value = 2941.986
$ws.Range("K34").Value = 1479.8572
$ws.Range("L34").Value = 2941.986
$ws.Range("M34").Value = -1277.8572
$ws.Range("N34").Value = -3345.986

$ws.Range("H58").Value = 2084
$ws.Range("I58").Value = 1448.8334
$ws.Range("K58").Value = 1448.8334
$ws.Range("M58").Value = -1245.8334

$ws.Range("H132").Value = 3031.5
$ws.Range("I132").Value = 2257.4443
$ws.Range("K132").Value = 6772.3329
$ws.Range("M132").Value = -4242.3329

$ws.Range("H134").Value = 2962.9333
$ws.Range("I134").Value = 1942.4615
$ws.Range("K134").Value = 5827.3845
$ws.Range("M134").Value = -3292.3845

$ws.Range("H136").Value = 2084
$ws.Range("I136").Value = 1448.8334
$ws.Range("K136").Value = 4346.5002
$ws.Range("M136").Value = -1796.5002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 295.6111
$ws.Range("J92").Value = 771.75
$ws.Range("L92").Value = 2315.25
$ws.Range("N92").Value = -4811.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10275.889
$ws.Range("I70").Value = 10104.714
$ws.Range("K70").Value = 10104.714
$ws.Range("M70").Value = -9834.714

$ws.Range("H73").Value = 10275.889
$ws.Range("I73").Value = 10104.714
$ws.Range("K73").Value = 10104.714
$ws.Range("M73").Value = -9168.714

$ws.Range("H80").Value = 2823.4443
$ws.Range("I80").Value = 2553.3333
$ws.Range("J80").Value = 2958.5
$ws.Range("K80").Value = 2553.3333
$ws.Range("L80").Value = 2958.5
$ws.Range("M80").Value = -1555.3333
$ws.Range("N80").Value = -4954.5

$ws.Range("H83").Value = 2823.4443
$ws.Range("I83").Value = 2553.3333
$ws.Range("J83").Value = 2958.5
$ws.Range("K83").Value = 12766.6665
$ws.Range("L83").Value = 14792.5
$ws.Range("M83").Value = -7774.666499999999
$ws.Range("N83").Value = -24776.5

$ws.Range("H113").Value = 4539.7
$ws.Range("I113").Value = 3000
$ws.Range("J113").Value = 4924.625
$ws.Range("K113").Value = 3000
$ws.Range("L113").Value = 4924.625
$ws.Range("M113").Value = -830
$ws.Range("N113").Value = -9264.625

$ws.Range("H122").Value = 3224.4211
$ws.Range("I122").Value = 2947.4285
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 8842.2855
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -6392.2855
$ws.Range("N122").Value = -16900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 425026.75
$ws.Range("I61").Value = 509134.34
$ws.Range("J61").Value = 4488.75
$ws.Range("K61").Value = 509134.34
$ws.Range("L61").Value = 4488.75
$ws.Range("M61").Value = -508932.34
$ws.Range("N61").Value = -4892.75

$ws.Range("H93").Value = 1250.15
$ws.Range("I93").Value = 1224.0555
$ws.Range("K93").Value = 1224.0555
$ws.Range("M93").Value = 23.94450000000006

$ws.Range("H113").Value = 425026.75
$ws.Range("I113").Value = 509134.34
$ws.Range("J113").Value = 4488.75
$ws.Range("K113").Value = 509134.34
$ws.Range("L113").Value = 4488.75
$ws.Range("M113").Value = -506964.34
$ws.Range("N113").Value = -8828.75

$ws.Range("H122").Value = 3587.6316
$ws.Range("I122").Value = 3090.077
$ws.Range("K122").Value = 9270.231
$ws.Range("M122").Value = -6820.231

$ws.Range("H136").Value = 6774.25
$ws.Range("I136").Value = 3611.5293
$ws.Range("K136").Value = 10834.5879
$ws.Range("M136").Value = -8284.5879

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 1631150.6
$ws.Range("I4").Value = 10001000
$ws.Range("J4").Value = 109359.91
$ws.Range("K4").Value = 10001000
$ws.Range("L4").Value = 109359.91
$ws.Range("M4").Value = -10000887
$ws.Range("N4").Value = -109585.91
